# Update "Max Cr" / "Max Date" values on the "Intermediate for Mapping" sheet
# and the matching "Max Cr" / "Max Date" values on the "Intermediate Exhibit"
# sheet. The source data stores every value as literal text (even numbers and
# dates), so we force the target ranges to Text format before writing the new
# values - otherwise COM Automation would silently reinterpret strings such as
# "7.61" or "2002-08-01" as a number / date and strip the original formatting.

$wb = $excel.ActiveWorkbook
$wsMap = $wb.Worksheets.Item("Intermediate for Mapping")
$wsEx  = $wb.Worksheets.Item("Intermediate Exhibit")

# Force the columns we are about to touch to store plain text so the literal
# strings below are preserved exactly (no auto-conversion to number/date).
$wsMap.Range("V2:V27").NumberFormat = "@"
$wsMap.Range("W2:W27").NumberFormat = "@"
$wsEx.Range("G4:G31").NumberFormat = "@"
$wsEx.Range("H4:H31").NumberFormat = "@"

# row on "Intermediate for Mapping" -> row on "Intermediate Exhibit"
$rows = @(
  @{ Map = 2;  Ex = 4;  V = "7.61";   W = "2002-08-01"; G = "7.61";   H = "8/1/02" },
  @{ Map = 9;  Ex = 12; V = "4.05";   W = "2010-06-17"; G = "4.05";   H = "6/17/10" },
  @{ Map = 11; Ex = 14; V = "3.84";   W = "2011-10-28"; G = "3.84";   H = "10/28/11" },
  @{ Map = 12; Ex = 16; V = "12.9";   W = "2008-01-24"; G = "12.9";   H = "1/24/08" },
  @{ Map = 13; Ex = 17; V = "2.67";   W = "2012-09-04"; G = "2.67";   H = "9/4/12" },
  @{ Map = 14; Ex = 18; V = "5.0";    W = $null;        G = "5";      H = $null },
  @{ Map = 15; Ex = 19; V = "5.6";    W = "2008-01-23"; G = "5.6";    H = "1/23/08" },
  @{ Map = 16; Ex = 20; V = "5.43";   W = "2011-03-10"; G = "5.43";   H = "3/10/11" },
  @{ Map = 17; Ex = 21; V = "4.7";    W = $null;        G = "4.7";    H = $null },
  @{ Map = 18; Ex = 22; V = "3.9425"; W = "2007-01-11"; G = "3.9425"; H = "1/11/07" },
  @{ Map = 21; Ex = 25; V = "9.63";   W = "2004-02-06"; G = "9.63";   H = "2/6/04" },
  @{ Map = 22; Ex = 26; V = "5.52";   W = "2010-04-29"; G = "5.52";   H = "4/29/10" },
  @{ Map = 24; Ex = 28; V = "2.6";    W = "2009-12-14"; G = "2.6";    H = "12/14/09" },
  @{ Map = 26; Ex = 30; V = "5.55";   W = "2009-12-04"; G = "5.55";   H = "12/4/09" },
  @{ Map = 27; Ex = 31; V = "5.6789"; W = "2010-06-17"; G = "5.6789"; H = "6/17/10" }
)

foreach ($row in $rows) {
    $wsMap.Cells.Item($row.Map, 22).Value = $row.V
    if ($row.W -ne $null) {
        $wsMap.Cells.Item($row.Map, 23).Value = $row.W
    }

    $wsEx.Cells.Item($row.Ex, 7).Value = $row.G
    if ($row.H -ne $null) {
        $wsEx.Cells.Item($row.Ex, 8).Value = $row.H
    }
}
